$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Row 48: Beta Build - MQTT communication work
$ws.Range("A48").Value = 46063
$ws.Range("B48").Value = "Beta Build"
$ws.Range("C48").Value = 0.44444444444444442
$ws.Range("D48").Value = 0.58333333333333337
$ws.Range("F48").Value = "Experimented with group and successfully got MQTT client to communicate t microcontroller. Started modifying main so the main loop would be triggered from MQTT message"

# Row 49: Beta Build - cleaning up MQTT branch
$ws.Range("A49").Value = 46063
$ws.Range("B49").Value = "Beta Build"
$ws.Range("C49").Value = 0.82013888888888886
$ws.Range("F49").Value = "Began cleaning up MQTT branch to prepare it for more clean merge"

# Update the view state to match where the user left off scrolling/selecting
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 3
$ws.Range("F53").Select()
